# Add two new test rows (16 and 17) for the built-in date formats
# "d-mmm" (numFmtId 16) and "d-mmm-yy" (numFmtId 15), which were not
# previously covered by this DateTime-detection test workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (15) as a formatting template for the
# two new rows, so font/fill/number-format base styling is carried over.
$ws.Range("A15:C15").Copy($ws.Range("A16:C16"))
$ws.Range("A15:C15").Copy($ws.Range("A17:C17"))

# Row 17 first: "Date 14-Mar" using the built-in d-mmm format (numFmtId 16)
$ws.Range("A17").Value2 = "Date 14-Mar"
$ws.Range("B17").Value2 = 40982.563138888887
$ws.Range("B17").NumberFormat = "d-mmm"
$ws.Range("C17").Formula = "=B17"

# Row 16 next: "Date 14-Mar-12" using the built-in d-mmm-yy format (numFmtId 15)
$ws.Range("A16").Value2 = "Date 14-Mar-12"
$ws.Range("B16").Value2 = 40982.563138888887
$ws.Range("B16").NumberFormat = "d-mmm-yy"
$ws.Range("C16").Formula = "=B16"

# Match the selection left behind by the author (active cell moved to A17)
[void]$ws.Range("A17").Select()
